$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 936.6923
$ws.Range("I39").Value = 68
$ws.Range("J39").Value = 1681.2858
$ws.Range("K39").Value = 204
$ws.Range("L39").Value = 5043.857400000001
$ws.Range("M39").Value = 92
$ws.Range("N39").Value = -5635.857400000001

# Row 132
$ws.Range("H132").Value = 3127404.5
$ws.Range("I132").Value = 3705948.2
$ws.Range("K132").Value = 11117844.6
$ws.Range("M132").Value = -11115314.6

# Row 137
$ws.Range("H137").Value = 3747.7073
$ws.Range("I137").Value = 3906.4062
$ws.Range("J137").Value = 3183.4443
$ws.Range("K137").Value = 11719.2186
$ws.Range("L137").Value = 9550.332900000001
$ws.Range("M137").Value = -9169.2186
$ws.Range("N137").Value = -14650.3329

# Row 138
$ws.Range("H138").Value = 5811.271
$ws.Range("I138").Value = 1688.7826
$ws.Range("J138").Value = 9603.959999999999
$ws.Range("K138").Value = 5066.3478
$ws.Range("L138").Value = 28811.88
$ws.Range("M138").Value = 73.65220000000045
$ws.Range("N138").Value = -39091.88

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 22730262
$ws.Range("I2").Value = 31251608
$ws.Range("J2").Value = 6670.6665
$ws.Range("K2").Value = 31251608
$ws.Range("L2").Value = 6670.6665
$ws.Range("M2").Value = -31251495
$ws.Range("N2").Value = -6896.6665

# Row 32
$ws.Range("H32").Value = 2595.36
$ws.Range("I32").Value = 2139.6345
$ws.Range("J32").Value = 8650
$ws.Range("K32").Value = 2139.6345
$ws.Range("L32").Value = 8650
$ws.Range("M32").Value = -1852.6345
$ws.Range("N32").Value = -9224

# Row 61
$ws.Range("H61").Value = 1388.075
$ws.Range("I61").Value = 729.2286
$ws.Range("K61").Value = 729.2286
$ws.Range("M61").Value = -517.2286

# Row 116
$ws.Range("H116").Value = 22730262
$ws.Range("I116").Value = 31251608
$ws.Range("J116").Value = 6670.6665
$ws.Range("K116").Value = 31251608
$ws.Range("L116").Value = 6670.6665
$ws.Range("M116").Value = -31249314
$ws.Range("N116").Value = -11258.6665

# Row 136
$ws.Range("H136").Value = 1388.075
$ws.Range("I136").Value = 729.2286
$ws.Range("K136").Value = 2187.6858
$ws.Range("M136").Value = 362.3141999999998

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 22730262
$ws.Range("I3").Value = 31251608
$ws.Range("J3").Value = 6670.6665
$ws.Range("K3").Value = 31251608
$ws.Range("L3").Value = 6670.6665
$ws.Range("M3").Value = -31251494
$ws.Range("N3").Value = -6898.6665

# Row 49
$ws.Range("H49").Value = 44032.5
$ws.Range("J49").Value = 44032.5
$ws.Range("L49").Value = 44032.5
$ws.Range("N49").Value = -44510.5

# Row 94
$ws.Range("H94").Value = 408.0909
$ws.Range("I94").Value = 422.2
$ws.Range("J94").Value = 377.85715
$ws.Range("K94").Value = 422.2
$ws.Range("L94").Value = 377.85715
$ws.Range("M94").Value = 28.80000000000001
$ws.Range("N94").Value = -1279.85715

# Row 99
$ws.Range("H99").Value = 2267.3572
$ws.Range("I99").Value = 1666.9524
$ws.Range("J99").Value = 4068.5715
$ws.Range("K99").Value = 1666.9524
$ws.Range("L99").Value = 4068.5715
$ws.Range("M99").Value = -168.9523999999999
$ws.Range("N99").Value = -7064.5715

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2834.7114
$ws.Range("I31").Value = 1630.9231
$ws.Range("J31").Value = 6446.077
$ws.Range("K31").Value = 1630.9231
$ws.Range("L31").Value = 6446.077
$ws.Range("M31").Value = -1335.9231
$ws.Range("N31").Value = -7036.077

# Row 34
$ws.Range("H34").Value = 2834.7114
$ws.Range("I34").Value = 1630.9231
$ws.Range("J34").Value = 6446.077
$ws.Range("K34").Value = 1630.9231
$ws.Range("L34").Value = 6446.077
$ws.Range("M34").Value = -1428.9231
$ws.Range("N34").Value = -6850.077

# Row 105
$ws.Range("H105").Value = 4002.353
$ws.Range("I105").Value = 3993.3333
$ws.Range("J105").Value = 4012.5
$ws.Range("K105").Value = 3993.3333
$ws.Range("L105").Value = 4012.5
$ws.Range("M105").Value = -2246.3333
$ws.Range("N105").Value = -7506.5

# Row 107
$ws.Range("H107").Value = 1446.08
$ws.Range("I107").Value = 408.17648
$ws.Range("J107").Value = 3651.625
$ws.Range("K107").Value = 408.17648
$ws.Range("L107").Value = 3651.625
$ws.Range("M107").Value = 1511.82352
$ws.Range("N107").Value = -7491.625

# Row 134
$ws.Range("H134").Value = 1732.1464
$ws.Range("I134").Value = 1119.037
$ws.Range("J134").Value = 2914.5715
$ws.Range("K134").Value = 3357.111
$ws.Range("L134").Value = 8743.7145
$ws.Range("M134").Value = -822.1109999999999
$ws.Range("N134").Value = -13813.7145

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4561.6553
$ws.Range("I70").Value = 4671.857
$ws.Range("J70").Value = 4272.375
$ws.Range("K70").Value = 4671.857
$ws.Range("L70").Value = 4272.375
$ws.Range("M70").Value = -4401.857
$ws.Range("N70").Value = -4812.375

# Row 73
$ws.Range("H73").Value = 4561.6553
$ws.Range("I73").Value = 4671.857
$ws.Range("J73").Value = 4272.375
$ws.Range("K73").Value = 4671.857
$ws.Range("L73").Value = 4272.375
$ws.Range("M73").Value = -3735.857
$ws.Range("N73").Value = -6144.375

# Row 80
$ws.Range("H80").Value = 2984.6924
$ws.Range("I80").Value = 2875
$ws.Range("J80").Value = 3160.2
$ws.Range("K80").Value = 2875
$ws.Range("L80").Value = 3160.2
$ws.Range("M80").Value = -1877
$ws.Range("N80").Value = -5156.2

# Row 83
$ws.Range("H83").Value = 2984.6924
$ws.Range("I83").Value = 2875
$ws.Range("J83").Value = 3160.2
$ws.Range("K83").Value = 14375
$ws.Range("L83").Value = 15801
$ws.Range("M83").Value = -9383
$ws.Range("N83").Value = -25785

# Row 113
$ws.Range("H113").Value = 1231.0952
$ws.Range("I113").Value = 842.8333
$ws.Range("J113").Value = 3560.6667
$ws.Range("K113").Value = 842.8333
$ws.Range("L113").Value = 3560.6667
$ws.Range("M113").Value = 1327.1667
$ws.Range("N113").Value = -7900.6667

# Row 132
$ws.Range("H132").Value = 3294.6978
$ws.Range("I132").Value = 3078.4666
$ws.Range("J132").Value = 3793.6924
$ws.Range("K132").Value = 9235.399800000001
$ws.Range("L132").Value = 11381.0772
$ws.Range("M132").Value = -6705.399800000001
$ws.Range("N132").Value = -16441.0772

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 90910750
$ws.Range("I61").Value = 200000700
$ws.Range("J61").Value = 2466.1667
$ws.Range("K61").Value = 200000700
$ws.Range("L61").Value = 2466.1667
$ws.Range("M61").Value = -200000498
$ws.Range("N61").Value = -2870.1667

# Row 113
$ws.Range("H113").Value = 90910750
$ws.Range("I113").Value = 200000700
$ws.Range("J113").Value = 2466.1667
$ws.Range("K113").Value = 200000700
$ws.Range("L113").Value = 2466.1667
$ws.Range("M113").Value = -199998530
$ws.Range("N113").Value = -6806.1667

$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 75014
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 75014
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 75014
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -75600
